# Auto-generated edit script: applies numeric cell updates described in the commit diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (35 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M98").Value = 769.6111
$ws.Range("H98").Value = 728.3889
$ws.Range("K98").Value = 728.3889
$ws.Range("I98").Value = 728.3889
$ws.Range("N120").Value = -71676
$ws.Range("H120").Value = 62000
$ws.Range("J120").Value = 62000
$ws.Range("L120").Value = 62000
$ws.Range("I122").Value = 728.3889
$ws.Range("H122").Value = 728.3889
$ws.Range("M122").Value = 264.8332999999998
$ws.Range("K122").Value = 2185.1667
$ws.Range("H129").Value = 3666.3333
$ws.Range("I129").Value = 2499.5
$ws.Range("K129").Value = 7498.5
$ws.Range("M129").Value = -2498.5
$ws.Range("K132").Value = 3470.8236
$ws.Range("M132").Value = -940.8235999999997
$ws.Range("I132").Value = 1156.9412
$ws.Range("H132").Value = 1195.3143
$ws.Range("I135").Value = 20001966
$ws.Range("M135").Value = -180015159
$ws.Range("K135").Value = 180017694
$ws.Range("H135").Value = 12505441
$ws.Range("M138").Value = -7988
$ws.Range("I138").Value = 4376
$ws.Range("L138").Value = 24352.941
$ws.Range("K138").Value = 13128
$ws.Range("H138").Value = 7267.273
$ws.Range("J138").Value = 8117.647
$ws.Range("N138").Value = -34632.941
$ws.Range("H141").Value = 5372.8184
$ws.Range("K141").Value = 16573.143
$ws.Range("M141").Value = -11393.143
$ws.Range("I141").Value = 5524.381

# --- Sheet: ARM (35 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M45").Value = -843
$ws.Range("K45").Value = 1220
$ws.Range("I45").Value = 1220
$ws.Range("H45").Value = 4394.143
$ws.Range("H61").Value = 4034.3257
$ws.Range("M61").Value = -2670.6667
$ws.Range("K61").Value = 2882.6667
$ws.Range("I61").Value = 2882.6667
$ws.Range("L101").Value = 55000
$ws.Range("N101").Value = -61490
$ws.Range("H101").Value = 55000
$ws.Range("J101").Value = 55000
$ws.Range("L103").Value = 79998.336
$ws.Range("N103").Value = -82342.336
$ws.Range("H103").Value = 79998.336
$ws.Range("J103").Value = 79998.336
$ws.Range("L105").Value = 93116.25
$ws.Range("N105").Value = -100104.25
$ws.Range("H105").Value = 93116.25
$ws.Range("J105").Value = 93116.25
$ws.Range("N122").Value = -23525.3638
$ws.Range("L122").Value = 18625.3638
$ws.Range("I122").Value = 2739.1177
$ws.Range("H122").Value = 4102.0713
$ws.Range("J122").Value = 6208.4546
$ws.Range("M122").Value = -5767.3531
$ws.Range("K122").Value = 8217.3531
$ws.Range("J131").Value = 66518.125
$ws.Range("H131").Value = 66518.125
$ws.Range("N131").Value = -76598.125
$ws.Range("L131").Value = 66518.125
$ws.Range("K136").Value = 8648.000100000001
$ws.Range("M136").Value = -6098.000100000001
$ws.Range("H136").Value = 4034.3257
$ws.Range("I136").Value = 2882.6667

# --- Sheet: BSM (4 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4600
$ws.Range("K22").Value = 5375
$ws.Range("M22").Value = -5202
$ws.Range("I22").Value = 5375

# --- Sheet: CRP (38 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5775.923
$ws.Range("M16").Value = -5921.9
$ws.Range("I16").Value = 6208.9
$ws.Range("K16").Value = 6208.9
$ws.Range("L31").Value = 4014.6155
$ws.Range("K31").Value = 2039.7826
$ws.Range("M31").Value = -1744.7826
$ws.Range("H31").Value = 2752.9167
$ws.Range("I31").Value = 2039.7826
$ws.Range("J31").Value = 4014.6155
$ws.Range("N31").Value = -4604.6155
$ws.Range("J34").Value = 4014.6155
$ws.Range("M34").Value = -1837.7826
$ws.Range("I34").Value = 2039.7826
$ws.Range("K34").Value = 2039.7826
$ws.Range("H34").Value = 2752.9167
$ws.Range("L34").Value = 4014.6155
$ws.Range("N34").Value = -4418.6155
$ws.Range("I58").Value = 1466.5834
$ws.Range("M58").Value = -1263.5834
$ws.Range("K58").Value = 1466.5834
$ws.Range("H58").Value = 2612.4707
$ws.Range("M113").Value = -4038.9
$ws.Range("H113").Value = 5775.923
$ws.Range("K113").Value = 6208.9
$ws.Range("I113").Value = 6208.9
$ws.Range("K132").Value = 8865.105599999999
$ws.Range("M132").Value = -6335.105599999999
$ws.Range("I132").Value = 2955.0352
$ws.Range("H132").Value = 3507.9436
$ws.Range("I134").Value = 2675.3333
$ws.Range("H134").Value = 4369.5454
$ws.Range("M134").Value = -5490.999899999999
$ws.Range("K134").Value = 8025.999899999999
$ws.Range("K136").Value = 4399.7502
$ws.Range("M136").Value = -1849.7502
$ws.Range("H136").Value = 2612.4707
$ws.Range("I136").Value = 1466.5834

# --- Sheet: CUL (4 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N122").Value = -32641.375
$ws.Range("L122").Value = 27741.375
$ws.Range("H122").Value = 2589.1538
$ws.Range("J122").Value = 3082.375

# --- Sheet: GSM (67 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("J70").Value = 5912.8
$ws.Range("K70").Value = 90207.69500000001
$ws.Range("M70").Value = -89937.69500000001
$ws.Range("I70").Value = 90207.69500000001
$ws.Range("L70").Value = 5912.8
$ws.Range("H70").Value = 66792.44500000001
$ws.Range("N70").Value = -6452.8
$ws.Range("M73").Value = -89271.69500000001
$ws.Range("H73").Value = 66792.44500000001
$ws.Range("I73").Value = 90207.69500000001
$ws.Range("K73").Value = 90207.69500000001
$ws.Range("N73").Value = -7784.8
$ws.Range("J73").Value = 5912.8
$ws.Range("L73").Value = 5912.8
$ws.Range("N76").ClearContents()
$ws.Range("L76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("I80").Value = 377370.34
$ws.Range("J80").Value = 500000000
$ws.Range("H80").Value = 125283020
$ws.Range("N80").Value = -500001996
$ws.Range("M80").Value = -376372.34
$ws.Range("L80").Value = 500000000
$ws.Range("K80").Value = 377370.34
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("M83").Value = -1881859.7
$ws.Range("I83").Value = 377370.34
$ws.Range("K83").Value = 1886851.7
$ws.Range("J83").Value = 500000000
$ws.Range("L83").Value = 2500000000
$ws.Range("H83").Value = 125283020
$ws.Range("N83").Value = -2500009984
$ws.Range("J85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("K97").Value = 2899.5
$ws.Range("H97").Value = 83335260
$ws.Range("I97").Value = 2899.5
$ws.Range("N97").Value = -250000992
$ws.Range("L97").Value = 250000000
$ws.Range("J97").Value = 250000000
$ws.Range("M97").Value = -2403.5
$ws.Range("L101").Value = 58000
$ws.Range("N101").Value = -64490
$ws.Range("H101").Value = 58000
$ws.Range("J101").Value = 58000
$ws.Range("K102").Value = 1772.5
$ws.Range("M102").Value = -150.5
$ws.Range("I102").Value = 1772.5
$ws.Range("H102").Value = 1612.2
$ws.Range("K132").Value = 35328
$ws.Range("M132").Value = -32798
$ws.Range("I132").Value = 11776
$ws.Range("H132").Value = 13648.647

# --- Sheet: LTW (16 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3149.7407
$ws.Range("N61").Value = -4895.857
$ws.Range("J61").Value = 4491.857
$ws.Range("L61").Value = 4491.857
$ws.Range("H93").Value = 5000
$ws.Range("L93").Value = 5001
$ws.Range("N93").Value = -7497
$ws.Range("J93").Value = 5001
$ws.Range("H113").Value = 3149.7407
$ws.Range("J113").Value = 4491.857
$ws.Range("L113").Value = 4491.857
$ws.Range("N113").Value = -8831.857
$ws.Range("K136").Value = 12448.2
$ws.Range("M136").Value = -9898.199999999999
$ws.Range("H136").Value = 6589.0713
$ws.Range("I136").Value = 4149.4

# --- Sheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M113").Value = 400.3000000000002
$ws.Range("H113").Value = 653.1070999999999
$ws.Range("K113").Value = 1769.7
$ws.Range("J113").Value = 811.125
$ws.Range("L113").Value = 2433.375
$ws.Range("N113").Value = -6773.375
$ws.Range("I113").Value = 589.9
$ws.Range("I122").Value = 1556.3684
$ws.Range("H122").Value = 2676.9565
$ws.Range("M122").Value = -2219.1052
$ws.Range("K122").Value = 4669.1052
$ws.Range("K136").Value = 10408.125
$ws.Range("M136").Value = -7858.125
$ws.Range("H136").Value = 4007.2058
$ws.Range("I136").Value = 3469.375

